$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing rows 2:34 down to 3:35
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the new first data point (month 0, df 1)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1
$ws.Range("B2").NumberFormat = "0.0000"

# Update selection to match the recorded state in the workbook
$ws.Range("B2").Select()
